# Rename the inline picture shapes that live in the document's headers and
# footers:
#   - footer 1 & footer 2: Pearson logo  image1.png -> image2.png
#   - header 2            : BTec logo    image2.jpg -> image1.jpg
#
# InlineShape has no settable "Name" of its own in the Word object model,
# so each picture is briefly converted to a floating Shape (which does
# expose .Name), renamed, then converted back to an inline picture so the
# layout/anchoring is left exactly as it was.

$d = $word.ActiveDocument
$sec = $d.Sections.Item(1)

function Rename-InlineLogo($range, $newName) {
    $inline = $range.InlineShapes.Item(1)
    $floating = $inline.ConvertToShape()
    $floating.Name = $newName
    [void]$floating.ConvertToInlineShape()
}

# Footer 1 - Pearson Edexcel logo
$footer1 = $sec.Footers.Item(1)
Rename-InlineLogo $footer1.Range "image2.png"

# Footer 2 - Pearson Edexcel logo
$footer2 = $sec.Footers.Item(2)
Rename-InlineLogo $footer2.Range "image2.png"

# Header 2 - BTec logo
$header2 = $sec.Headers.Item(2)
Rename-InlineLogo $header2.Range "image1.jpg"
